$wb = $excel.ActiveWorkbook

# --- Sheet "Kupci": replace customer #1 values ---
$wsKupci = $wb.Worksheets.Item("Kupci")
$wsKupci.Range("C2").Value = "Kompanija1 d.o.o."
$wsKupci.Range("C3").Value = "Kompanija2 d.d"
$wsKupci.Range("B2").Value = "21111111114"
$wsKupci.Range("B3").Value = "81111111110"
$wsKupci.Range("B3").Select()

# --- Sheet "Racuni": update oib references for the invoices ---
$wsRacuni = $wb.Worksheets.Item("Racuni")
$wsRacuni.Range("G2").Value = "21111111114"
$wsRacuni.Range("G3").Value = "21111111114"
$wsRacuni.Range("G4").Value = "81111111110"
$wsRacuni.Range("C10").Select()
